$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dtypes")

# Move the "bus_geodata" rows (currently rows 218-220) down to rows 230-232 so
# the new "pwl_cost" / "poly_cost" dtype tables can be placed ahead of them.
# Using Copy (instead of Rows.Insert) keeps cell formatting/styles identical to
# the source instead of synthesizing a brand-new style entry.
$ws.Range("A218:D220").Copy($ws.Range("A230"))
$excel.CutCopyMode = $false

# Reuse the formatting (bold + border, style used across column A) for the new
# id cells that will hold the "pwl_cost" / "poly_cost" rows.
$ws.Range("A217").Copy($ws.Range("A218:A229"))
$excel.CutCopyMode = $false

# New data rows describing the "pwl_cost" and "poly_cost" tables.
$rows = @(
    @(216, "pwl_cost",  "power_type",        "object"),
    @(217, "pwl_cost",  "element",           "uint32"),
    @(218, "pwl_cost",  "et",                "object"),
    @(219, "pwl_cost",  "points",            "object"),
    @(220, "poly_cost", "element",           "uint32"),
    @(221, "poly_cost", "et",                "object"),
    @(222, "poly_cost", "cp0_eur",           "float64"),
    @(223, "poly_cost", "cp1_eur_per_mw",    "float64"),
    @(224, "poly_cost", "cp2_eur_per_mw2",   "float64"),
    @(225, "poly_cost", "cq0_eur",           "float64"),
    @(226, "poly_cost", "cq1_eur_per_mvar",  "float64"),
    @(227, "poly_cost", "cq2_eur_per_mvar2", "float64")
)

$r = 218
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Renumber the id column for the "bus_geodata" rows that were pushed down
# by the insertion (they keep their data, only the running id changes).
$ws.Cells.Item(230, 1).Value = 228
$ws.Cells.Item(231, 1).Value = 229
$ws.Cells.Item(232, 1).Value = 230

$ws.Range("T238").Select() | Out-Null
